# Fruta / hortaliza, semanal
# Insert a new weekly price-report row for "Femacal de La Calera - Espinaca"
# at row 347 (pushing the existing rows 347-377 down to 348-378).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 347..377 down to 348..378 to make room for the new record.
$ws.Rows(347).Insert()

# Populate the newly inserted row with the new market observation.
$ws.Cells.Item(347, 1).Value  = 3
$ws.Cells.Item(347, 2).Value  = "Femacal de La Calera"
$ws.Cells.Item(347, 3).Value  = "Coquimbo"
$ws.Cells.Item(347, 4).Value  = 44769
$ws.Cells.Item(347, 5).Value  = 5
$ws.Cells.Item(347, 6).Value  = 100112012
$ws.Cells.Item(347, 7).Value  = "Espinaca"
$ws.Cells.Item(347, 8).Value  = "Sin especificar"
$ws.Cells.Item(347, 9).Value  = "Primera"
$ws.Cells.Item(347, 10).Value = 80
$ws.Cells.Item(347, 11).Value = 4000
$ws.Cells.Item(347, 12).Value = 4000
$ws.Cells.Item(347, 13).Value = 4000
$ws.Cells.Item(347, 14).Value = "$/docena de atados (3 kilos)"
$ws.Cells.Item(347, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(347, 16).Value = 1333
$ws.Cells.Item(347, 17).Value = 3
$ws.Cells.Item(347, 18).Value = "Hortaliza"
